$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Write the literal text via a formula (so Excel's literal-to-date /
    # literal-to-number auto-detection never kicks in), then convert the
    # formula back down to a plain value in place, and finally strip any
    # style that got picked up from the column defaults so the cell matches
    # its neighbours (no explicit style attribute).
    $escaped = $text -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $ws.Range($addr).Borders.LineStyle = 0
}

# Row 59 used to be a lone placeholder row holding only the "**" marker in
# column A. Turn it into a full equipment-log entry and push the "**"
# marker down to the new row 60.
Set-TextCell "A59" "05/03/2018"
Set-TextCell "B59" "3070"
$ws.Range("C59").Value = 3
$ws.Range("C59").Borders.LineStyle = 0
Set-TextCell "D59" "test"
Set-TextCell "E59" "1"
Set-TextCell "F59" "80000571"
Set-TextCell "G59" "246"

Set-TextCell "A60" "**"

# Excel re-registers the sheet's _FilterDatabase defined name whenever the
# filtered range is touched/re-saved; replicate that extra bookkeeping
# entry.
$wb.Names.Add("_xlnm._FilterDatabase", "=main!`$A`$2:`$G`$2")
